# Update currency year from 2023 to 2024 (and related 2012-base ratio / notes)
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Update the descriptive labels / text that reference the currency year
$about.Range("A24").Value = "2024 dollars"
$about.Range("A26").Value = 0.73191600598044548
$about.Range("B26").Value = "2024 dollars per 2012 dollar"
$about.Range("B29").Value = 'which in this case is "2012 dollars per 2024 dollar."'
$about.Range("B30").Value = "2012 dollars are worth more than 2024 dollars, so we need a"

# Move the active selection to B31 to match the saved workbook state
$about.Range("B31").Select()
